# Add a new "Decision Tree" result row (Fuller dataset, no demographics,
# no hyperparameter tuning) to Sheet1, then tidy up the AutoFilter/
# _FilterDatabase bookkeeping to reflect the new data extent.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Append the new row of data (row 22) -----------------------------
$ws.Range("A22").Value = "Decision Tree"
$ws.Range("B22").Value = $false
$ws.Range("C22").Value = $false
$ws.Range("D22").Value = 0.97
$ws.Range("E22").Value = 0.78
$ws.Range("F22").Value = 0.54
$ws.Range("G22").Value = 0.64
$ws.Range("H22").Value = 0.89
$ws.Range("I22").Value = "Fuller dataset"

# Match the percentage number format used by the rest of column D
# (copy formatting only from the row above, which already uses it).
$ws.Range("D21").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Move the active selection to the new last cell -------------------
$null = $ws.Range("I22").Select()

# --- Remove the AutoFilter (it no longer covers the new row) ----------
if ($ws.AutoFilterMode) {
    $ws.AutoFilterMode = $false
}

# --- Update the hidden _FilterDatabase defined name to the new extent -
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$I`$21"
    }
}
